$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (7 -> 8) ---
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667

# --- Data updates for rows 2-5 ---
$ws.Range("A2").Value = 45111.50694444445
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("A3").Value = 45111.51388888889
$ws.Range("B3").Value = 19.734
$ws.Range("C3").Value = 14.753
$ws.Range("D3").Value = 0.588
$ws.Range("E3").Value = 43.398
$ws.Range("F3").Value = 37.041
$ws.Range("G3").Value = 15.615
$ws.Range("H3").Value = 51.47
$ws.Range("I3").Value = 23.895
$ws.Range("J3").Value = 10.895
$ws.Range("K3").Value = 16.684
$ws.Range("L3").Value = 18.933
$ws.Range("M3").Value = 18.299
$ws.Range("N3").Value = 5.999
$ws.Range("O3").Value = 15.651
$ws.Range("P3").Value = 22.36
$ws.Range("Q3").Value = 12.757
$ws.Range("R3").Value = 0.297
$ws.Range("S3").Value = 0.918
$ws.Range("T3").Value = 232.427
$ws.Range("U3").Value = 43.203
$ws.Range("V3").Value = 14.422
$ws.Range("W3").Value = 29.377
$ws.Range("X3").Value = 15.258
$ws.Range("Y3").Value = 2.985
$ws.Range("Z3").Value = 26.585
$ws.Range("AA3").Value = 12.795
$ws.Range("AB3").Value = 11.434
$ws.Range("AC3").Value = 13.462
$ws.Range("AD3").Value = 19.514
$ws.Range("AE3").Value = 1.059
$ws.Range("AF3").Value = 47.493
$ws.Range("AG3").Value = 8.731999999999999
$ws.Range("AH3").Value = 17.789
$ws.Range("A4").Value = 45111.52083333334
$ws.Range("B4").Value = 16.547
$ws.Range("C4").Value = 12.385
$ws.Range("D4").Value = 0.442
$ws.Range("E4").Value = 36.377
$ws.Range("F4").Value = 30.611
$ws.Range("G4").Value = 13.108
$ws.Range("H4").Value = 51.269
$ws.Range("I4").Value = 20.042
$ws.Range("J4").Value = 9.134
$ws.Range("K4").Value = 13.794
$ws.Range("L4").Value = 15.374
$ws.Range("M4").Value = 15.333
$ws.Range("N4").Value = 4.775
$ws.Range("O4").Value = 13.026
$ws.Range("P4").Value = 18.753
$ws.Range("Q4").Value = 10.747
$ws.Range("R4").Value = 0.196
$ws.Range("S4").Value = 0.711
$ws.Range("T4").Value = 192.665
$ws.Range("U4").Value = 36.347
$ws.Range("V4").Value = 12.014
$ws.Range("W4").Value = 24.705
$ws.Range("X4").Value = 12.825
$ws.Range("Y4").Value = 2.267
$ws.Range("Z4").Value = 25.059
$ws.Range("AA4").Value = 10.701
$ws.Range("AB4").Value = 9.552
$ws.Range("AC4").Value = 11.239
$ws.Range("AD4").Value = 16.004
$ws.Range("AE4").Value = 0.664
$ws.Range("AF4").Value = 47.327
$ws.Range("AG4").Value = 7.186
$ws.Range("AH4").Value = 14.923
$ws.Range("A5").Value = 45111.52777777778
$ws.Range("B5").Value = 24.31
$ws.Range("C5").Value = 18.23
$ws.Range("D5").Value = 0.71
$ws.Range("E5").Value = 53.19
$ws.Range("F5").Value = 44.31
$ws.Range("G5").Value = 19.21
$ws.Range("H5").Value = 73.65000000000001
$ws.Range("I5").Value = 29.45
$ws.Range("J5").Value = 13.33
$ws.Range("K5").Value = 19.98
$ws.Range("L5").Value = 21.85
$ws.Range("M5").Value = 22.48
$ws.Range("N5").Value = 6.55
$ws.Range("O5").Value = 19.06
$ws.Range("P5").Value = 27.39
$ws.Range("Q5").Value = 15.79
$ws.Range("R5").Value = 0.24
$ws.Range("S5").Value = 0.89
$ws.Range("T5").Value = 284.33
$ws.Range("U5").Value = 53.25
$ws.Range("V5").Value = 17.59
$ws.Range("W5").Value = 36.15
$ws.Range("X5").Value = 18.88
$ws.Range("Y5").Value = 2.9
$ws.Range("Z5").Value = 36.05
$ws.Range("AA5").Value = 15.62
$ws.Range("AB5").Value = 13.83
$ws.Range("AC5").Value = 16.28
$ws.Range("AD5").Value = 22.89
$ws.Range("AE5").Value = 0.49
$ws.Range("AF5").Value = 67.27
$ws.Range("AG5").Value = 10.28
$ws.Range("AH5").Value = 21.94

# --- Remove the old row 6 (data now only spans rows 1-5) ---
$ws.Rows.Item(6).Delete()
